$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 716.5263
$ws.Range("I19").Value = 610.4545000000001
$ws.Range("J19").Value = 862.375
$ws.Range("K19").Value = 610.4545000000001
$ws.Range("L19").Value = 862.375
$ws.Range("M19").Value = -435.4545000000001
$ws.Range("N19").Value = -1212.375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 906.5599999999999
$ws.Range("I33").Value = 219.3158
$ws.Range("J33").Value = 3082.8333
$ws.Range("K33").Value = 219.3158
$ws.Range("L33").Value = 3082.8333
$ws.Range("M33").Value = 9.684200000000004
$ws.Range("N33").Value = -3540.8333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 664.7222
$ws.Range("I92").Value = 497.66666
$ws.Range("J92").Value = 1500
$ws.Range("K92").Value = 497.66666
$ws.Range("L92").Value = 1500
$ws.Range("M92").Value = 750.33334
$ws.Range("N92").Value = -3996

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4906789.5
$ws.Range("I132").Value = 5560694.5
$ws.Range("K132").Value = 16682083.5
$ws.Range("M132").Value = -16679553.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 59477.234
$ws.Range("I2").Value = 719
$ws.Range("J2").Value = 143417.58
$ws.Range("K2").Value = 719
$ws.Range("L2").Value = 143417.58
$ws.Range("M2").Value = -606
$ws.Range("N2").Value = -143643.58

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 101443.4
$ws.Range("I45").Value = 126415.375
$ws.Range("J45").Value = 1555.5
$ws.Range("K45").Value = 126415.375
$ws.Range("L45").Value = 1555.5
$ws.Range("M45").Value = -126038.375
$ws.Range("N45").Value = -2309.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 59477.234
$ws.Range("I116").Value = 719
$ws.Range("J116").Value = 143417.58
$ws.Range("K116").Value = 719
$ws.Range("L116").Value = 143417.58
$ws.Range("M116").Value = 1575
$ws.Range("N116").Value = -148005.58

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1924.4546
$ws.Range("I122").Value = 1781.9375
$ws.Range("J122").Value = 2304.5
$ws.Range("K122").Value = 5345.8125
$ws.Range("L122").Value = 6913.5
$ws.Range("M122").Value = -2895.8125
$ws.Range("N122").Value = -11813.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3328.0312
$ws.Range("I132").Value = 3500.9167
$ws.Range("J132").Value = 2809.375
$ws.Range("K132").Value = 10502.7501
$ws.Range("L132").Value = 8428.125
$ws.Range("M132").Value = -7972.750100000001
$ws.Range("N132").Value = -13488.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 59477.234
$ws.Range("I3").Value = 719
$ws.Range("J3").Value = 143417.58
$ws.Range("K3").Value = 719
$ws.Range("L3").Value = 143417.58
$ws.Range("M3").Value = -605
$ws.Range("N3").Value = -143645.58

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1107.7778
$ws.Range("J80").Value = 1291
$ws.Range("L80").Value = 1291
$ws.Range("N80").Value = -3287

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 1107.7778
$ws.Range("J83").Value = 1291
$ws.Range("L83").Value = 6455
$ws.Range("N83").Value = -16439

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 66698470
$ws.Range("I107").Value = 125057660
$ws.Range("J107").Value = 2247.5715
$ws.Range("K107").Value = 125057660
$ws.Range("L107").Value = 2247.5715
$ws.Range("M107").Value = -125055740
$ws.Range("N107").Value = -6087.5715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 5509.5
$ws.Range("J29").Value = 10000
$ws.Range("L29").Value = 10000
$ws.Range("N29").Value = -10586

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 24833.484
$ws.Range("I31").Value = 1406.7097
$ws.Range("K31").Value = 1406.7097
$ws.Range("M31").Value = -1111.7097

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 24833.484
$ws.Range("I34").Value = 1406.7097
$ws.Range("K34").Value = 1406.7097
$ws.Range("M34").Value = -1204.7097

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 32980
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 32980
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 32980
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -35270

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 565.65515
$ws.Range("I107").Value = 577.2632
$ws.Range("J107").Value = 543.6
$ws.Range("K107").Value = 577.2632
$ws.Range("L107").Value = 543.6
$ws.Range("M107").Value = 1342.7368
$ws.Range("N107").Value = -4383.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2133.5715
$ws.Range("I122").Value = 846
$ws.Range("J122").Value = 2648.6
$ws.Range("K122").Value = 2538
$ws.Range("L122").Value = 7945.799999999999
$ws.Range("M122").Value = -88
$ws.Range("N122").Value = -12845.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1206.1915
$ws.Range("I5").Value = 1027.9375
$ws.Range("J5").Value = 1298.1936
$ws.Range("K5").Value = 3083.8125
$ws.Range("L5").Value = 3894.5808
$ws.Range("M5").Value = -2971.8125
$ws.Range("N5").Value = -4118.5808

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 71.55556
$ws.Range("I38").Value = 46.666668
$ws.Range("J38").Value = 84
$ws.Range("K38").Value = 140.000004
$ws.Range("L38").Value = 252
$ws.Range("M38").Value = 206.999996
$ws.Range("N38").Value = -946

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 2855.5557
$ws.Range("I58").Value = 3000
$ws.Range("J58").Value = 2837.5
$ws.Range("K58").Value = 9000
$ws.Range("L58").Value = 8512.5
$ws.Range("M58").Value = -8872
$ws.Range("N58").Value = -8768.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1206.1915
$ws.Range("I135").Value = 1027.9375
$ws.Range("J135").Value = 1298.1936
$ws.Range("K135").Value = 9251.4375
$ws.Range("L135").Value = 11683.7424
$ws.Range("M135").Value = -6716.4375
$ws.Range("N135").Value = -16753.7424

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 10709.833
$ws.Range("I138").Value = 15779.714
$ws.Range("J138").Value = 3612
$ws.Range("K138").Value = 47339.142
$ws.Range("L138").Value = 10836
$ws.Range("M138").Value = -42199.142
$ws.Range("N138").Value = -21116

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 10710.909
$ws.Range("I141").Value = 11182
$ws.Range("K141").Value = 33546
$ws.Range("M141").Value = -28366

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 9133.333000000001
$ws.Range("J26").Value = 9133.333000000001
$ws.Range("L26").Value = 9133.333000000001
$ws.Range("N26").Value = -9693.333000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H50").Value = 9133.333000000001
$ws.Range("J50").Value = 9133.333000000001
$ws.Range("L50").Value = 9133.333000000001
$ws.Range("N50").Value = -10129.333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1933.1666
$ws.Range("I102").Value = 1200.16
$ws.Range("J102").Value = 3599.0908
$ws.Range("K102").Value = 1200.16
$ws.Range("L102").Value = 3599.0908
$ws.Range("M102").Value = 421.8399999999999
$ws.Range("N102").Value = -6843.0908

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 874.65515
$ws.Range("I122").Value = 690.3077
$ws.Range("J122").Value = 1024.4375
$ws.Range("K122").Value = 2070.9231
$ws.Range("L122").Value = 3073.3125
$ws.Range("M122").Value = 379.0769
$ws.Range("N122").Value = -7973.3125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2180.5625
$ws.Range("I126").Value = 2268.375
$ws.Range("J126").Value = 2092.75
$ws.Range("K126").Value = 6805.125
$ws.Range("L126").Value = 6278.25
$ws.Range("M126").Value = -4335.125
$ws.Range("N126").Value = -11218.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2975.6287
$ws.Range("I132").Value = 2725.7856
$ws.Range("J132").Value = 3975
$ws.Range("K132").Value = 8177.3568
$ws.Range("L132").Value = 11925
$ws.Range("M132").Value = -5647.3568
$ws.Range("N132").Value = -16985

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4243.12
$ws.Range("I132").Value = 4026.3635
$ws.Range("J132").Value = 5832.6665
$ws.Range("K132").Value = 12079.0905
$ws.Range("L132").Value = 17497.9995
$ws.Range("M132").Value = -9549.0905
$ws.Range("N132").Value = -22557.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 33469.5
$ws.Range("J56").Value = 43292.668
$ws.Range("L56").Value = 43292.668
$ws.Range("N56").Value = -44720.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1497.9231
$ws.Range("I126").Value = 1307.4
$ws.Range("K126").Value = 3922.2
$ws.Range("M126").Value = -1452.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2801.7354
$ws.Range("I132").Value = 2653.8276
$ws.Range("K132").Value = 7961.4828
$ws.Range("M132").Value = -5431.4828
